{"js": "// Replace the 25 \"a\u00d7b=c\" answer cells in the multiplication table with\n// their updated values, matching the XML diff 1:1 (unique old text ->\n// unique new text, applied in document order).\nconst replacements = [\n  [\"44\u00d776=3344\", \"67\u00d765=4355\"],\n  [\"65\u00d718=1170\", \"76\u00d760=4560\"],\n  [\"80\u00d756=4480\", \"84\u00d770=5880\"],\n  [\"70\u00d763=4410\", \"18\u00d744=792\"],\n  [\"77\u00d748=3696\", \"71\u00d743=3053\"],\n  [\"67\u00d727=1809\", \"80\u00d791=7280\"],\n  [\"82\u00d728=2296\", \"60\u00d797=5820\"],\n  [\"28\u00d739=1092\", \"82\u00d782=6724\"],\n  [\"55\u00d768=3740\", \"79\u00d782=6478\"],\n  [\"91\u00d784=7644\", \"87\u00d743=3741\"],\n  [\"44\u00d747=2068\", \"23\u00d732=736\"],\n  [\"68\u00d754=3672\", \"83\u00d787=7221\"],\n  [\"83\u00d753=4399\", \"69\u00d758=4002\"],\n  [\"98\u00d746=4508\", \"83\u00d776=6308\"],\n  [\"79\u00d721=1659\", \"11\u00d720=220\"],\n  [\"29\u00d794=2726\", \"21\u00d714=294\"],\n  [\"22\u00d749=1078\", \"11\u00d752=572\"],\n  [\"39\u00d766=2574\", \"31\u00d722=682\"],\n  [\"65\u00d722=1430\", \"51\u00d746=2346\"],\n  [\"82\u00d740=3280\", \"57\u00d756=3192\"],\n  [\"47\u00d747=2209\", \"64\u00d729=1856\"],\n  [\"27\u00d719=513\", \"41\u00d752=2132\"],\n  [\"36\u00d797=3492\", \"97\u00d715=1455\"],\n  [\"53\u00d721=1113\", \"47\u00d798=4606\"],\n  [\"26\u00d779=2054\", \"93\u00d711=1023\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"a\u00d7b=c\" answer cells in the multiplication table with\n# their updated values, matching the XML diff 1:1 (unique old text ->\n# unique new text, applied via Find/Replace over the whole document body).\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($old, $new) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n\nReplace-Text \"44\u00d776=3344\" \"67\u00d765=4355\"\nReplace-Text \"65\u00d718=1170\" \"76\u00d760=4560\"\nReplace-Text \"80\u00d756=4480\" \"84\u00d770=5880\"\nReplace-Text \"70\u00d763=4410\" \"18\u00d744=792\"\nReplace-Text \"77\u00d748=3696\" \"71\u00d743=3053\"\nReplace-Text \"67\u00d727=1809\" \"80\u00d791=7280\"\nReplace-Text \"82\u00d728=2296\" \"60\u00d797=5820\"\nReplace-Text \"28\u00d739=1092\" \"82\u00d782=6724\"\nReplace-Text \"55\u00d768=3740\" \"79\u00d782=6478\"\nReplace-Text \"91\u00d784=7644\" \"87\u00d743=3741\"\nReplace-Text \"44\u00d747=2068\" \"23\u00d732=736\"\nReplace-Text \"68\u00d754=3672\" \"83\u00d787=7221\"\nReplace-Text \"83\u00d753=4399\" \"69\u00d758=4002\"\nReplace-Text \"98\u00d746=4508\" \"83\u00d776=6308\"\nReplace-Text \"79\u00d721=1659\" \"11\u00d720=220\"\nReplace-Text \"29\u00d794=2726\" \"21\u00d714=294\"\nReplace-Text \"22\u00d749=1078\" \"11\u00d752=572\"\nReplace-Text \"39\u00d766=2574\" \"31\u00d722=682\"\nReplace-Text \"65\u00d722=1430\" \"51\u00d746=2346\"\nReplace-Text \"82\u00d740=3280\" \"57\u00d756=3192\"\nReplace-Text \"47\u00d747=2209\" \"64\u00d729=1856\"\nReplace-Text \"27\u00d719=513\" \"41\u00d752=2132\"\nReplace-Text \"36\u00d797=3492\" \"97\u00d715=1455\"\nReplace-Text \"53\u00d721=1113\" \"47\u00d798=4606\"\nReplace-Text \"26\u00d779=2054\" \"93\u00d711=1023\"\n"}
